$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3 (Invalid) changes from 0 to 1
$ws.Range("G3").Value = 1

# H3:H18 (Absent) all change from 0 to 1
$ws.Range("H3:H18").Value = 1
